# Update salt prices and loadings in all scenarios
# (Magnesium chloride unit price row 8, Zinc sulfate unit price row 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Magnesium chloride unit price: loading + lower/upper bounds
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 - Zinc sulfate unit price: loading + lower/upper bounds
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Reflect the active selection left by the author after editing rows 8-9
$ws.Range("A8:XFD9").Select()
